$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD, AE, AF) reusing the existing bold/
# centered/bordered header style from AC1, then set their captions.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate season-record columns for every data row (2-65) with this
# team's Wins/Losses/Ties totals.
$ws.Range("AD2:AD65").Value = 90
$ws.Range("AE2:AE65").Value = 72
$ws.Range("AF2:AF65").Value = 0
